# Updated symbol list with GitHub Actions: refresh Price (D) and
# Volume(1h) (E) columns with the latest scraped coinranking.com figures.
# Values are written as literal text (leading "'" quote-prefix forces
# text entry so strings like "288.00" / "1.07%" are preserved verbatim
# instead of being coerced into numbers), then the cell style is reset
# back to "Normal" so no residual quote-prefix formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.07%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'29.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.36%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.201"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'5.28%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06960"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'4.16%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.440"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.55%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.556"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'5.17%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.403"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'3.47%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9041"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-3.53%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1601"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.63%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07719"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'16.76%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07720"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.24%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.02930"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.21%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09023"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.26%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001603"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.71%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0006501"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.92%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006249"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.79%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.473"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.52%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.29%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3234"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.94%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1339"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.39%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.016"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.23%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'4.80%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'1.02%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001210"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.65%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004149"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.41%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001169"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-6.30%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'3.27%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04391"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.82%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006932"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.76%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'-0.13%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002068"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'2.58%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-4.28%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005800"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.83%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-1.86%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.34%"
$ws.Range("E47").Style = "Normal"
